$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPbES")

# Rename the existing "hydrogen" entry to "hydrogen combustion turbine"
$ws.Range("A24").Value = "hydrogen combustion turbine"

# Add the new row for "hydrogen combined cycle"
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25:AK25").Value = 3

# Apply the pasted-in look (black font + vertically centered) to A24 first,
# then copy that exact formatting onto A25 in a single atomic operation so
# no extra intermediate cell style is left behind in styles.xml.
$ws.Range("A24").Font.Color = 0
$ws.Range("A24").VerticalAlignment = -4108
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A25").PasteSpecial(-4122) | Out-Null

# Update the view: select the new row's data range.
$ws.Activate() | Out-Null
$ws.Range("B25:AK25").Select() | Out-Null

# Restore the originally active sheet so DPbES doesn't end up as the
# selected tab.
$wb.Worksheets.Item("About").Activate() | Out-Null
